# The whole dataset (Timestamp column A and the matching "Lookup" text in
# column E) is shifted forward by exactly one day - row-for-row the same
# Import/Export/Quarter values stay put, only the dates move on by 1 day
# (15/16/17 Feb 2026 -> 16/17/18 Feb 2026). This mirrors adding a fresh
# day's worth of rows onto the rolling Entsoe "Unintended Deviation" feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 195

for ($r = 2; $r -le $lastRow; $r++) {

    # --- Column A: Timestamp (numeric date serial) -> + 1 day ---
    $cellA = $ws.Cells.Item($r, 1)
    $serial = $cellA.Value2()
    $cellA.Value = $serial + 1

    # --- Column E: Lookup ("DD.MM.YYYY" + quarter-suffix) -> + 1 day ---
    $cellE = $ws.Cells.Item($r, 5)
    $lookup = $cellE.Value2()

    $day = [int]$lookup.Substring(0, 2)
    $month = $lookup.Substring(3, 2)
    $yearAndSuffix = $lookup.Substring(6)

    $newDay = $day + 1
    $newLookup = "{0:D2}.{1}.{2}" -f $newDay, $month, $yearAndSuffix

    $cellE.Value = $newLookup
}
